# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "datos actualizados" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 11:34"

# 2. Plain data refreshes (country stays in the same row, only counts change)
# Row 6  - India
$ws.Range("B6").Value = 1339176
$ws.Range("C6").Value = 2154
$ws.Range("D6").Value = 850303
$ws.Range("E6").Value = 457448

# Row 20 - Banglades
$ws.Range("B20").Value = 221178
$ws.Range("C20").Value = 2520
$ws.Range("D20").Value = 122090
$ws.Range("E20").Value = 96214
$ws.Range("G20").Value = 38
$ws.Range("H20").Value = 2874

# Row 27 - Indonesia
$ws.Range("B27").Value = 97286
$ws.Range("C27").Value = 1868
$ws.Range("D27").Value = 55354
$ws.Range("E27").Value = 37218
$ws.Range("G27").Value = 49
$ws.Range("H27").Value = 4714

# Row 33 - Filipinas
$ws.Range("B33").Value = 78412
$ws.Range("C33").Value = 1968
$ws.Range("D33").Value = 25752
$ws.Range("E33").Value = 50763
$ws.Range("G33").Value = 18
$ws.Range("H33").Value = 1897

# Row 34 - Oman
$ws.Range("B34").Value = 74858
$ws.Range("C34").Value = 1067
$ws.Range("D34").Value = 54061
$ws.Range("E34").Value = 20426
$ws.Range("G34").Value = 12
$ws.Range("H34").Value = 371

# Row 54 - Afganistan
$ws.Range("E54").Value = 9995
$ws.Range("G54").Value = 23
$ws.Range("H54").Value = 1248

# Row 64 - Austria
$ws.Range("B64").Value = 20338
$ws.Range("C64").Value = 124
$ws.Range("D64").Value = 18124
$ws.Range("E64").Value = 1502
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 712

# Row 86 - Malasia
$ws.Range("B86").Value = 8884
$ws.Range("C86").Value = 23
$ws.Range("D86").Value = 8594
$ws.Range("E86").Value = 167

# Row 89 - Finlandia
$ws.Range("B89").Value = 7388
$ws.Range("C89").Value = 8
$ws.Range("E89").Value = 139

# Row 114 - Sri Lanka
$ws.Range("D114").Value = 2103
$ws.Range("E114").Value = 650

# Row 119 - Hong Kong
$ws.Range("E119").Value = 949
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 17

# 3. Reorder countries: Singapur now appears before Portugal (rows 45-46),
#    with Singapur's stats refreshed and Portugal's unchanged.
$ws.Range("A45").Value = "Singapur"
$ws.Range("B45").Value = 49888
$ws.Range("C45").Value = 513
$ws.Range("D45").Value = 45172
$ws.Range("E45").Value = 4689
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 27

$ws.Range("A46").Value = "Portugal"
$ws.Range("B46").Value = 49692
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 34687
$ws.Range("E46").Value = 13293
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 1712

# 4. Reorder countries: Polonia now appears before Rumania (rows 48-49),
#    with Polonia's stats refreshed and Rumania's unchanged.
$ws.Range("A48").Value = "Polonia"
$ws.Range("B48").Value = 42622
$ws.Range("C48").Value = 584
$ws.Range("D48").Value = 31997
$ws.Range("E48").Value = 8961
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 9
$ws.Range("H48").Value = 1664

$ws.Range("A49").Value = "Rumania"
$ws.Range("B49").Value = 42394
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 25349
$ws.Range("E49").Value = 14895
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 2150

# 5. Swap Islas Malvinas and Groenlandia (rows 210-211); underlying data unchanged.
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"
